$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 11 (the existing "76442780" Cash payment) into row 12 for the
# new payment event, then adjust the new row's timestamp.
$ws.Range("A11:G11").Copy($ws.Range("A12:G12"))
$ws.Range("D12").Value = "2025-08-15T10:00:30"

# Normalize A11's phone number to a real number (it was stored as text).
$ws.Range("A11").Value = 76442780
